# fix: datetime timezone issues
#
# Add a "creation" column (E) to the "client" sheet holding a fixed,
# literal timestamp string ("2024-07-22 12:34:56"). The value is entered
# as quote-prefixed text (leading apostrophe) so Excel stores it verbatim
# as a string instead of converting it to a serial date/time value, which
# is what caused the timezone-related discrepancies this commit fixes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("client")

# Header for the new column
$ws.Range("E1").Value = "creation"

# Literal timestamp, forced to text via leading apostrophe (quotePrefix)
$ws.Range("E2").Value = "'2024-07-22 12:34:56"

# Widen the new column to fit its content
$ws.Columns.Item(5).ColumnWidth = 17.3

# Mirror the author's final selection/active cell
$ws.Range("E2").Select()
